# Insert 3 new leaf columns (region, marketplace_id, marketplace_name)
# into Raw_One right before the existing seller_origin column (E), shifting
# seller_origin.. and everything after it 3 columns to the right.

$wb = $excel.ActiveWorkbook
$rawOne = $wb.Worksheets.Item("Raw_One")

# xlShiftToRight = -4161
$rawOne.Range("E1:G1").Insert(-4161)

$rawOne.Range("E1").Value = "region"
$rawOne.Range("F1").Value = "marketplace_id"
$rawOne.Range("G1").Value = "marketplace_name"

# Re-point every defined name that sat at/after the old column E so it
# tracks the same header cell at its new (shifted +3) column. Range.Insert
# does not retarget named ranges on its own, so this has to be explicit.
$wb.Names.Item("seller_origin").RefersTo            = "=Raw_One!`$H`$1:`$H`$1"
$wb.Names.Item("seller_origin_level1").RefersTo     = "=Raw_One!`$I`$1:`$I`$1"
$wb.Names.Item("seller_origin_level2").RefersTo     = "=Raw_One!`$J`$1:`$J`$1"
$wb.Names.Item("wtd_threep_net_ord_gms").RefersTo   = "=Raw_One!`$K`$1:`$K`$1"
$wb.Names.Item("wtd_fba_net_ord_gms").RefersTo      = "=Raw_One!`$L`$1:`$L`$1"
$wb.Names.Item("wtd_mfn_net_ord_gms").RefersTo      = "=Raw_One!`$M`$1:`$M`$1"
$wb.Names.Item("wtd_threep_net_ord_units").RefersTo = "=Raw_One!`$N`$1:`$N`$1"
$wb.Names.Item("wtd_fba_net_ord_units").RefersTo    = "=Raw_One!`$O`$1:`$O`$1"
$wb.Names.Item("wtd_mfn_net_ord_units").RefersTo    = "=Raw_One!`$P`$1:`$P`$1"
$wb.Names.Item("mtd_threep_net_ord_gms").RefersTo   = "=Raw_One!`$Q`$1:`$Q`$1"
$wb.Names.Item("mtd_fba_net_ord_gms").RefersTo      = "=Raw_One!`$R`$1:`$R`$1"
$wb.Names.Item("mtd_mfn_net_ord_gms").RefersTo      = "=Raw_One!`$S`$1:`$S`$1"
$wb.Names.Item("mtd_threep_net_ord_units").RefersTo = "=Raw_One!`$T`$1:`$T`$1"
$wb.Names.Item("mtd_fba_net_ord_units").RefersTo    = "=Raw_One!`$U`$1:`$U`$1"
$wb.Names.Item("mtd_mfn_net_ord_units").RefersTo    = "=Raw_One!`$V`$1:`$V`$1"
$wb.Names.Item("qtd_threep_net_ord_gms").RefersTo   = "=Raw_One!`$W`$1:`$W`$1"
$wb.Names.Item("qtd_fba_net_ord_gms").RefersTo      = "=Raw_One!`$X`$1:`$X`$1"
$wb.Names.Item("qtd_mfn_net_ord_gms").RefersTo      = "=Raw_One!`$Y`$1:`$Y`$1"
$wb.Names.Item("qtd_threep_net_ord_units").RefersTo = "=Raw_One!`$Z`$1:`$Z`$1"
$wb.Names.Item("qtd_fba_net_ord_units").RefersTo    = "=Raw_One!`$AA`$1:`$AA`$1"
$wb.Names.Item("qtd_mfn_net_ord_units").RefersTo    = "=Raw_One!`$AB`$1:`$AB`$1"
$wb.Names.Item("ytd_threep_net_ord_gms").RefersTo   = "=Raw_One!`$AC`$1:`$AC`$1"
$wb.Names.Item("ytd_fba_net_ord_gms").RefersTo      = "=Raw_One!`$AD`$1:`$AD`$1"
$wb.Names.Item("ytd_mfn_net_ord_gms").RefersTo      = "=Raw_One!`$AE`$1:`$AE`$1"
$wb.Names.Item("ytd_threep_net_ord_units").RefersTo = "=Raw_One!`$AF`$1:`$AF`$1"
$wb.Names.Item("ytd_fba_net_ord_units").RefersTo    = "=Raw_One!`$AG`$1:`$AG`$1"
$wb.Names.Item("ytd_mfn_net_ord_units").RefersTo    = "=Raw_One!`$AH`$1:`$AH`$1"

# Brand-new defined names for the 3 newly inserted columns.
$wb.Names.Add("region", "=Raw_One!`$E`$1:`$E`$1")
$wb.Names.Add("marketplace_id", "=Raw_One!`$F`$1:`$F`$1")
$wb.Names.Add("marketplace_name", "=Raw_One!`$G`$1:`$G`$1")

# --- Main sheet: relabel the metric rows with the wtd_ prefix and add the
# SUMIFS roll-ups in J:O now that "region" is a real column/named range. ---
$main = $wb.Worksheets.Item("Main")

$regions = @("NA", "EU", "JP", "Emerging")
$metricPairs = @(
    @("fba_net_ord_gms",   "wtd_fba_net_ord_gms",   12),
    @("fba_net_ord_units", "wtd_fba_net_ord_units", 16),
    @("mfn_net_ord_gms",   "wtd_mfn_net_ord_gms",   20),
    @("mfn_net_ord_units", "wtd_mfn_net_ord_units", 24),
    @("net_ord_gms",       "wtd_net_ord_gms",       28),
    @("net_ord_units",     "wtd_net_ord_units",     32)
)

$cols = @("J", "K", "L", "M", "N", "O")

foreach ($pair in $metricPairs) {
    $newLabel = $pair[1]
    $startRow = $pair[2]
    for ($i = 0; $i -lt $regions.Length; $i++) {
        $r = $startRow + $i
        $region = $regions[$i]
        $main.Range("A" + $r).Value = $newLabel
        foreach ($col in $cols) {
            $cell = $main.Range($col + $r)
            $cell.Formula = "=SUMIFS(" + $newLabel + ",activity_year," + $col + "4,activity_week," + $col + "3,region,`"" + $region + "`")"
        }
    }
}
